$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.955.78'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.902.20'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -4.56%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.47'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.92%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3817'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.55'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07704'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -5.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9808'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.01'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.001.80'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.968'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.669'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.07033'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.006'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '83.92'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -5.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000009531'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -5.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.75'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -4.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '28.950.26'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.324'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.89'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.160.79'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.102'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '157.50'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.583'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -6.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '117.54'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.850'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09288'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8596'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -6.19%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.250'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -7.74%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.018'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05691'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.146'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.004'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.460'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -5.64%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '9.273'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -6.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.744'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5189'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.96%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -6.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.085'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06822'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '111.24'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.52%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -5.30%  '
